$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.361.11"
$ws.Range("E2").Value = "  -2.31%  "
$ws.Range("D3").Value = "'3.473.96"
$ws.Range("E3").Value = "  -4.53%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'577.29"
$ws.Range("E5").Value = "  -4.56%  "
$ws.Range("D6").Value = "'192.66"
$ws.Range("E6").Value = "  -3.67%  "
$ws.Range("D7").Value = "'0.609"
$ws.Range("E7").Value = "  -3.06%  "
$ws.Range("D8").Value = "'3.463.35"
$ws.Range("E8").Value = "  -4.55%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'0.205"
$ws.Range("E10").Value = "  -7.30%  "
$ws.Range("D11").Value = "'0.617"
$ws.Range("D12").Value = "'51.42"
$ws.Range("E12").Value = "  -4.72%  "
$ws.Range("D13").Value = "'0.0000285"
$ws.Range("E13").Value = "  -7.14%  "
$ws.Range("D14").Value = "'9.13"
$ws.Range("E14").Value = "  -4.68%  "
$ws.Range("D15").Value = "'4.036.50"
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("D16").Value = "'644.28"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "'69.146.35"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").Value = "'3.463.93"
$ws.Range("E18").Value = "  -5.49%  "
$ws.Range("D19").Value = "'12.33"
$ws.Range("E19").Value = "  -5.21%  "
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").Value = "'18.17"
$ws.Range("E21").Value = "  -4.69%  "
$ws.Range("E22").Value = "  -5.85%  "
$ws.Range("D23").Value = "'17.82"
$ws.Range("E23").Value = "  -4.22%  "
$ws.Range("D24").Value = "'5.29"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").Value = "'99.11"
$ws.Range("E25").Value = "  -5.00%  "
$ws.Range("D26").Value = "'4.28"
$ws.Range("E26").Value = "  -7.73%  "
$ws.Range("E27").Value = "  -4.77%  "
$ws.Range("D28").Value = "'9.92"
$ws.Range("E28").Value = "  -5.25%  "
$ws.Range("D29").Value = "'9.32"
$ws.Range("E29").Value = "  -4.48%  "
$ws.Range("D30").Value = "'32.42"
$ws.Range("E30").Value = "  -4.54%  "
$ws.Range("D31").Value = "'4.27"
$ws.Range("E31").Value = "  -11.25%  "
$ws.Range("E32").Value = "  -6.63%  "
$ws.Range("D33").Value = "'11.60"
$ws.Range("E33").Value = "  -5.35%  "
$ws.Range("E34").Value = "  -5.59%  "
$ws.Range("D35").Value = "'60.88"
$ws.Range("D36").Value = "'3.712.02"
$ws.Range("E36").Value = "  -8.02%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'523.40"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "'0.0₃0792"
$ws.Range("E39").Value = "  -10.24%  "
$ws.Range("D40").Value = "'2.94"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("E42").Value = "  -4.88%  "
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("B44").Value = "CoreDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D44").Value = "'3.50"
$ws.Range("E44").Value = "  +69.00%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'34.25"
$ws.Range("E45").Value = "  -6.80%  "
$ws.Range("E46").Value = "  -4.38%  "
$ws.Range("D47").Value = "'3.36"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").Value = "'2.82"
$ws.Range("E48").Value = "  -6.54%  "
$ws.Range("E49").Value = "  -4.67%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").Value = "'8.14"
$ws.Range("E51").Value = "  -6.20%  "
